$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.234.23"
$ws.Range("D3").Value = "1.857.75"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'0.7086"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").Value = "'238.41"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.07978"
$ws.Range("E8").Value = "  +3.48%  "
$ws.Range("D9").Value = "'0.3029"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'23.45"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").Value = "'0.08201"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").Value = "'5.178"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7029"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'89.62"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.739.19"
$ws.Range("E15").Value = "  -6.35%  "
$ws.Range("D16").Value = "29.096.87"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'5.810"
$ws.Range("D18").Value = "'0.000007883"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").Value = "'13.24"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "'237.63"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "'0.9991"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.033.03"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'7.468"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").Value = "'162.88"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").Value = "'8.902"
$ws.Range("D27").Value = "'0.1438"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "'18.10"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "'1.916"
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("D30").Value = "'1.422"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'1.477"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'4.360"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").Value = "'4.019"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "'0.05186"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "'1.158"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "'0.9958"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "'0.01850"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'2.721"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("D41").Value = "'0.9315"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "1.139.30"
$ws.Range("E42").Value = "  +5.49%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.908"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4251"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'70.15"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'102.35"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").Value = "'0.5330"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("D49").Value = "'1.762"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "'9.164"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'6.948"
$ws.Range("E51").Value = "  -0.82%  "
